$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $cell.Formula = '="' + $val + '"'
    $cell.Copy() | Out-Null
    $cell.PasteSpecial(-4163) | Out-Null
}

$ws.Cells.Item(2, 1).Value = "12:00"
$ws.Cells.Item(2, 2).Value = "12:25"
$ws.Cells.Item(2, 3).Value = "Table 1"
$ws.Cells.Item(2, 4).Value = "Group phase"
$ws.Cells.Item(2, 5).Value = "Group A"
$ws.Cells.Item(2, 6).Value = "Bálint Bakos"
$ws.Cells.Item(2, 7).Value = "2*"
$ws.Cells.Item(2, 8).Value = "Ákos Szvetnik"
Set-TextValue ($ws.Cells.Item(2, 9)) "0"
$ws.Cells.Item(3, 1).Value = "12:00"
$ws.Cells.Item(3, 2).Value = "12:25"
$ws.Cells.Item(3, 3).Value = "Table 2"
$ws.Cells.Item(3, 4).Value = "Group phase"
$ws.Cells.Item(3, 5).Value = "Group A"
$ws.Cells.Item(3, 6).Value = "Dani Pusztai"
Set-TextValue ($ws.Cells.Item(3, 7)) "1"
$ws.Cells.Item(3, 8).Value = "Máté Vass"
$ws.Cells.Item(3, 9).Value = "2*"
$ws.Cells.Item(4, 1).Value = "12:30"
$ws.Cells.Item(4, 2).Value = "12:55"
$ws.Cells.Item(4, 3).Value = "Table 1"
$ws.Cells.Item(4, 4).Value = "Group phase"
$ws.Cells.Item(4, 5).Value = "Group A"
$ws.Cells.Item(4, 6).Value = "Dáni Gyulai-Nagy"
Set-TextValue ($ws.Cells.Item(4, 7)) "0"
$ws.Cells.Item(4, 8).Value = "Laci Ferenczi"
$ws.Cells.Item(4, 9).Value = "2*"
$ws.Cells.Item(5, 1).Value = "12:30"
$ws.Cells.Item(5, 2).Value = "12:55"
$ws.Cells.Item(5, 3).Value = "Table 2"
$ws.Cells.Item(5, 4).Value = "Group phase"
$ws.Cells.Item(5, 5).Value = "Group A"
$ws.Cells.Item(5, 6).Value = "Berci Pusztai"
Set-TextValue ($ws.Cells.Item(5, 7)) "0"
$ws.Cells.Item(5, 8).Value = "Ákos Szvetnik"
$ws.Cells.Item(5, 9).Value = "2*"
$ws.Cells.Item(6, 1).Value = "13:00"
$ws.Cells.Item(6, 2).Value = "13:25"
$ws.Cells.Item(6, 3).Value = "Table 1"
$ws.Cells.Item(6, 4).Value = "Group phase"
$ws.Cells.Item(6, 5).Value = "Group A"
$ws.Cells.Item(6, 6).Value = "Bálint Bakos"
Set-TextValue ($ws.Cells.Item(6, 7)) "1"
$ws.Cells.Item(6, 8).Value = "Laci Ferenczi"
$ws.Cells.Item(6, 9).Value = "2*"
$ws.Cells.Item(7, 1).Value = "13:00"
$ws.Cells.Item(7, 2).Value = "13:25"
$ws.Cells.Item(7, 3).Value = "Table 2"
$ws.Cells.Item(7, 4).Value = "Group phase"
$ws.Cells.Item(7, 5).Value = "Group A"
$ws.Cells.Item(7, 6).Value = "Dani Pusztai"
Set-TextValue ($ws.Cells.Item(7, 7)) "0"
$ws.Cells.Item(7, 8).Value = "Dáni Gyulai-Nagy"
$ws.Cells.Item(7, 9).Value = "2*"
$ws.Cells.Item(8, 1).Value = "13:30"
$ws.Cells.Item(8, 2).Value = "13:55"
$ws.Cells.Item(8, 3).Value = "Table 1"
$ws.Cells.Item(8, 4).Value = "Group phase"
$ws.Cells.Item(8, 5).Value = "Group A"
$ws.Cells.Item(8, 6).Value = "Berci Pusztai"
Set-TextValue ($ws.Cells.Item(8, 7)) "1"
$ws.Cells.Item(8, 8).Value = "Máté Vass"
$ws.Cells.Item(8, 9).Value = "2*"
$ws.Cells.Item(9, 1).Value = "13:30"
$ws.Cells.Item(9, 2).Value = "13:55"
$ws.Cells.Item(9, 3).Value = "Table 2"
$ws.Cells.Item(9, 4).Value = "Group phase"
$ws.Cells.Item(9, 5).Value = "Group A"
$ws.Cells.Item(9, 6).Value = "Ákos Szvetnik"
Set-TextValue ($ws.Cells.Item(9, 7)) "0"
$ws.Cells.Item(9, 8).Value = "Laci Ferenczi"
$ws.Cells.Item(9, 9).Value = "2*"
$ws.Cells.Item(10, 1).Value = "14:00"
$ws.Cells.Item(10, 2).Value = "14:25"
$ws.Cells.Item(10, 3).Value = "Table 1"
$ws.Cells.Item(10, 4).Value = "Group phase"
$ws.Cells.Item(10, 5).Value = "Group A"
$ws.Cells.Item(10, 6).Value = "Bálint Bakos"
$ws.Cells.Item(10, 7).Value = "2*"
$ws.Cells.Item(10, 8).Value = "Dani Pusztai"
Set-TextValue ($ws.Cells.Item(10, 9)) "0"
$ws.Cells.Item(11, 1).Value = "14:00"
$ws.Cells.Item(11, 2).Value = "14:25"
$ws.Cells.Item(11, 3).Value = "Table 2"
$ws.Cells.Item(11, 4).Value = "Group phase"
$ws.Cells.Item(11, 5).Value = "Group A"
$ws.Cells.Item(11, 6).Value = "Berci Pusztai"
$ws.Cells.Item(11, 7).Value = "2*"
$ws.Cells.Item(11, 8).Value = "Laci Ferenczi"
Set-TextValue ($ws.Cells.Item(11, 9)) "0"
$ws.Cells.Item(12, 1).Value = "14:30"
$ws.Cells.Item(12, 2).Value = "14:55"
$ws.Cells.Item(12, 3).Value = "Table 1"
$ws.Cells.Item(12, 4).Value = "Group phase"
$ws.Cells.Item(12, 5).Value = "Group A"
$ws.Cells.Item(12, 6).Value = "Máté Vass"
$ws.Cells.Item(12, 7).Value = "2*"
$ws.Cells.Item(12, 8).Value = "Dáni Gyulai-Nagy"
Set-TextValue ($ws.Cells.Item(12, 9)) "1"
$ws.Cells.Item(13, 1).Value = "14:30"
$ws.Cells.Item(13, 2).Value = "14:55"
$ws.Cells.Item(13, 3).Value = "Table 2"
$ws.Cells.Item(13, 4).Value = "Group phase"
$ws.Cells.Item(13, 5).Value = "Group A"
$ws.Cells.Item(13, 6).Value = "Ákos Szvetnik"
$ws.Cells.Item(13, 7).Value = "2*"
$ws.Cells.Item(13, 8).Value = "Dani Pusztai"
Set-TextValue ($ws.Cells.Item(13, 9)) "0"
$ws.Cells.Item(14, 1).Value = "15:00"
$ws.Cells.Item(14, 2).Value = "15:25"
$ws.Cells.Item(14, 3).Value = "Table 1"
$ws.Cells.Item(14, 4).Value = "Group phase"
$ws.Cells.Item(14, 5).Value = "Group A"
$ws.Cells.Item(14, 6).Value = "Berci Pusztai"
$ws.Cells.Item(14, 7).Value = "2*"
$ws.Cells.Item(14, 8).Value = "Dáni Gyulai-Nagy"
Set-TextValue ($ws.Cells.Item(14, 9)) "0"
$ws.Cells.Item(15, 1).Value = "15:00"
$ws.Cells.Item(15, 2).Value = "15:25"
$ws.Cells.Item(15, 3).Value = "Table 2"
$ws.Cells.Item(15, 4).Value = "Group phase"
$ws.Cells.Item(15, 5).Value = "Group A"
$ws.Cells.Item(15, 6).Value = "Laci Ferenczi"
$ws.Cells.Item(15, 7).Value = "2*"
$ws.Cells.Item(15, 8).Value = "Dani Pusztai"
Set-TextValue ($ws.Cells.Item(15, 9)) "1"
$ws.Cells.Item(16, 1).Value = "15:30"
$ws.Cells.Item(16, 2).Value = "15:55"
$ws.Cells.Item(16, 3).Value = "Table 1"
$ws.Cells.Item(16, 4).Value = "Group phase"
$ws.Cells.Item(16, 5).Value = "Group A"
$ws.Cells.Item(16, 6).Value = "Máté Vass"
$ws.Cells.Item(16, 7).Value = "2*"
$ws.Cells.Item(16, 8).Value = "Bálint Bakos"
Set-TextValue ($ws.Cells.Item(16, 9)) "0"
$ws.Cells.Item(17, 1).Value = "15:30"
$ws.Cells.Item(17, 2).Value = "15:55"
$ws.Cells.Item(17, 3).Value = "Table 2"
$ws.Cells.Item(17, 4).Value = "Group phase"
$ws.Cells.Item(17, 5).Value = "Group A"
$ws.Cells.Item(17, 6).Value = "Berci Pusztai"
$ws.Cells.Item(17, 7).Value = "2*"
$ws.Cells.Item(17, 8).Value = "Dani Pusztai"
Set-TextValue ($ws.Cells.Item(17, 9)) "1"
$ws.Cells.Item(18, 1).Value = "16:00"
$ws.Cells.Item(18, 2).Value = "16:25"
$ws.Cells.Item(18, 3).Value = "Table 1"
$ws.Cells.Item(18, 4).Value = "Group phase"
$ws.Cells.Item(18, 5).Value = "Group A"
$ws.Cells.Item(18, 6).Value = "Dáni Gyulai-Nagy"
$ws.Cells.Item(18, 7).Value = "2*"
$ws.Cells.Item(18, 8).Value = "Bálint Bakos"
Set-TextValue ($ws.Cells.Item(18, 9)) "1"
$ws.Cells.Item(19, 1).Value = "16:00"
$ws.Cells.Item(19, 2).Value = "16:25"
$ws.Cells.Item(19, 3).Value = "Table 2"
$ws.Cells.Item(19, 4).Value = "Group phase"
$ws.Cells.Item(19, 5).Value = "Group A"
$ws.Cells.Item(19, 6).Value = "Máté Vass"
Set-TextValue ($ws.Cells.Item(19, 7)) "1"
$ws.Cells.Item(19, 8).Value = "Ákos Szvetnik"
$ws.Cells.Item(19, 9).Value = "2*"
$ws.Cells.Item(20, 1).Value = "16:30"
$ws.Cells.Item(20, 2).Value = "16:55"
$ws.Cells.Item(20, 3).Value = "Table 1"
$ws.Cells.Item(20, 4).Value = "Group phase"
$ws.Cells.Item(20, 5).Value = "Group A"
$ws.Cells.Item(20, 6).Value = "Berci Pusztai"
$ws.Cells.Item(20, 7).Value = "2*"
$ws.Cells.Item(20, 8).Value = "Bálint Bakos"
Set-TextValue ($ws.Cells.Item(20, 9)) "0"
$ws.Cells.Item(21, 1).Value = "16:30"
$ws.Cells.Item(21, 2).Value = "16:55"
$ws.Cells.Item(21, 3).Value = "Table 2"
$ws.Cells.Item(21, 4).Value = "Group phase"
$ws.Cells.Item(21, 5).Value = "Group A"
$ws.Cells.Item(21, 6).Value = "Dáni Gyulai-Nagy"
$ws.Cells.Item(21, 7).Value = "2*"
$ws.Cells.Item(21, 8).Value = "Ákos Szvetnik"
Set-TextValue ($ws.Cells.Item(21, 9)) "1"
$ws.Cells.Item(22, 1).Value = "17:00"
$ws.Cells.Item(22, 2).Value = "17:25"
$ws.Cells.Item(22, 3).Value = "Table 1"
$ws.Cells.Item(22, 4).Value = "Group phase"
$ws.Cells.Item(22, 5).Value = "Group A"
$ws.Cells.Item(22, 6).Value = "Laci Ferenczi"
$ws.Cells.Item(22, 7).Value = "2*"
$ws.Cells.Item(22, 8).Value = "Máté Vass"
Set-TextValue ($ws.Cells.Item(22, 9)) "1"
$ws.Cells.Item(23, 1).Value = "17:30"
$ws.Cells.Item(23, 2).Value = "17:55"
$ws.Cells.Item(23, 3).Value = "Table 1"
$ws.Cells.Item(23, 4).Value = "Knockout phase"
$ws.Cells.Item(23, 5).Value = "Match B1"
$ws.Cells.Item(23, 6).Value = "Laci Ferenczi"
$ws.Cells.Item(23, 7).Value = "3*"
$ws.Cells.Item(23, 8).Value = "Ákos Szvetnik"
Set-TextValue ($ws.Cells.Item(23, 9)) "1"
$ws.Cells.Item(24, 1).Value = "17:45"
$ws.Cells.Item(24, 2).Value = "18:10"
$ws.Cells.Item(24, 3).Value = "Table 2"
$ws.Cells.Item(24, 4).Value = "Knockout phase"
$ws.Cells.Item(24, 5).Value = "Match B2"
$ws.Cells.Item(24, 6).Value = "Berci Pusztai"
$ws.Cells.Item(24, 7).Value = "3*"
$ws.Cells.Item(24, 8).Value = "Máté Vass"
Set-TextValue ($ws.Cells.Item(24, 9)) "2"
$ws.Cells.Item(25, 1).Value = "17:30"
$ws.Cells.Item(25, 2).Value = "17:55"
$ws.Cells.Item(25, 3).Value = "Table 1"
$ws.Cells.Item(25, 4).Value = "Knockout phase"
$ws.Cells.Item(25, 5).Value = "Match B3"
$ws.Cells.Item(25, 6).Value = "Laci Ferenczi"
Set-TextValue ($ws.Cells.Item(25, 7)) "2"
$ws.Cells.Item(25, 8).Value = "Berci Pusztai"
$ws.Cells.Item(25, 9).Value = "4*"

$excel.CutCopyMode = 0
